$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constants (xlEdgeBottom = 9, xlContinuous = 1, xlMedium = -4138, xlColorIndexAutomatic = 1)
$xlEdgeBottom = 9
$xlMedium = -4138
$xlAutomatic = 1

function Set-BottomMediumBorder($cell) {
    $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlMedium
    $cell.Borders.Item($xlEdgeBottom).Weight = $xlMedium
    $cell.Borders.Item($xlEdgeBottom).ColorIndex = $xlAutomatic
}

# Row 2: L2 gains the same bottom border as the rest of the header-separator row
Set-BottomMediumBorder($ws.Cells.Item(2, 12))

# Row 3: L3 becomes the "2020" column header (bold 9pt Times New Roman + bottom border)
$L3 = $ws.Cells.Item(3, 12)
$L3.Font.Name = "Times New Roman"
$L3.Font.Size = 9
$L3.Font.Bold = $true
Set-BottomMediumBorder($L3)
$L3.Value = 2020

# Row 4: L4 new data value (regular 9pt Times New Roman, one-decimal number format)
$L4 = $ws.Cells.Item(4, 12)
$L4.Font.Name = "Times New Roman"
$L4.Font.Size = 9
$L4.Font.Bold = $false
$L4.NumberFormat = "0.0"
$L4.Value = 10

# Row 5: L5 new data value
$L5 = $ws.Cells.Item(5, 12)
$L5.Font.Name = "Times New Roman"
$L5.Font.Size = 9
$L5.Font.Bold = $false
$L5.Value = 6.4

# Row 6: L6 new data value
$L6 = $ws.Cells.Item(6, 12)
$L6.Font.Name = "Times New Roman"
$L6.Font.Size = 9
$L6.Font.Bold = $false
$L6.Value = 13.5

# Row 7: L7 new data value
$L7 = $ws.Cells.Item(7, 12)
$L7.Font.Name = "Times New Roman"
$L7.Font.Size = 9
$L7.Font.Bold = $false
$L7.Value = 24.3

# Row 8: L8 new data value
$L8 = $ws.Cells.Item(8, 12)
$L8.Font.Name = "Times New Roman"
$L8.Font.Size = 9
$L8.Font.Bold = $false
$L8.Value = 27.8

# Row 9: L9 new data value (previously no cell at all in this row)
$L9 = $ws.Cells.Item(9, 12)
$L9.Font.Name = "Times New Roman"
$L9.Font.Size = 9
$L9.Font.Bold = $false
$L9.Value = 20.9

# Row 10: L10 new data value
$L10 = $ws.Cells.Item(10, 12)
$L10.Font.Name = "Times New Roman"
$L10.Font.Size = 9
$L10.Font.Bold = $false
$L10.Value = 26.7

# Row 11: L11 new data value
$L11 = $ws.Cells.Item(11, 12)
$L11.Font.Name = "Times New Roman"
$L11.Font.Size = 9
$L11.Font.Bold = $false
$L11.Value = 28.4

# Row 12: L12 new data value (bottom row, gets the table's bottom border too)
$L12 = $ws.Cells.Item(12, 12)
$L12.Font.Name = "Times New Roman"
$L12.Font.Size = 9
$L12.Font.Bold = $false
Set-BottomMediumBorder($L12)
$L12.Value = 25

# Move/record the active selection as it was left after editing (M16)
$ws.Range("M16").Select()
